$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = "sku"
    $ws.Cells.Item($r, 2).Value = "name"
    $ws.Cells.Item($r, 3).Value = "quantity"
    $ws.Cells.Item($r, 4).Value = "cost_per"
    $ws.Cells.Item($r, 5).Value = "total_cost"
}
